$wb = $excel.ActiveWorkbook

# --- 1. Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- 3. Copy formatting from "Weekly Quantity" so the new sheet matches
#        the workbook's existing look (bold/centered/bordered header row,
#        date-formatted first column) without introducing new style entries.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A46").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- 4. Write headers ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- 5. Write data rows (ds, PO_Forecast, yhat_lower, yhat_upper) ---
$data = @(
    @(45067.99999999999, 71.0, -52.1469913282838, 199.3218361780529),
    @(45095.99999999999, 71.0, -61.37489699582044, 201.173608202817),
    @(45123.99999999999, 71.0, -47.53876503461261, 201.819355600353),
    @(45137.99999999999, 71.0, -52.69700016248021, 194.0677200787829),
    @(45144.99999999999, 71.0, -52.47642415048195, 200.5715861864387),
    @(45151.99999999999, 71.0, -61.59947894748497, 199.7976701564387),
    @(45158.99999999999, 71.0, -49.41135927828356, 187.7007438686732),
    @(45165.99999999999, 71.0, -51.32046383689503, 195.4761143823089),
    @(45179.99999999999, 71.0, -52.48384143295105, 200.8193742532641),
    @(45200.99999999999, 71.0, -66.7548486808327, 193.2275804401924),
    @(45214.99999999999, 71.0, -53.04113966103117, 198.3829682791275),
    @(45221.99999999999, 71.0, -51.85916455991468, 191.934194319674),
    @(45235.99999999999, 71.0, -55.13729965066985, 204.7210959289454),
    @(45242.99999999999, 71.0, -49.35244410790973, 193.4545569577919),
    @(45249.99999999999, 71.0, -50.63447311151499, 194.540484585443),
    @(45256.99999999999, 71.0, -57.4264124814884, 197.8516071278385),
    @(45263.99999999999, 71.0, -54.26366716344275, 196.650222542691),
    @(45270.99999999999, 71.0, -56.12662845267637, 194.1211452046792),
    @(45277.99999999999, 71.0, -55.61643059134654, 199.3411224029486),
    @(45298.99999999999, 71.0, -46.92231587706723, 195.1017834775699),
    @(45445.99999999999, 72.0, -56.74353398837473, 206.948240800566),
    @(45452.99999999999, 72.0, -53.28253787831981, 186.5219320881281),
    @(45459.99999999999, 72.0, -52.50963880578345, 190.1923678325104),
    @(45466.99999999999, 72.0, -54.77817100362211, 195.2889184750318),
    @(45487.99999999999, 72.0, -55.24039483245295, 192.2188152633357),
    @(45508.99999999999, 72.0, -53.15749846783901, 199.3711050094181),
    @(45529.99999999999, 72.0, -52.07314232595023, 201.8085975321499),
    @(45536.99999999999, 72.0, -43.61861453051986, 200.8582484940482),
    @(45564.99999999999, 72.0, -52.30537393987356, 208.7379845346554),
    @(45578.99999999999, 72.0, -51.99541558984632, 195.8886250103356),
    @(45585.99999999999, 72.0, -62.69073908282876, 197.2212938438475),
    @(45592.99999999999, 72.0, -58.85952243824094, 193.9033356410292),
    @(45599.99999999999, 72.0, -60.29802452818003, 196.2373804921115),
    @(45613.99999999999, 72.0, -52.60712939346978, 200.1875017463815),
    @(45620.99999999999, 72.0, -44.04400814060485, 201.8753514337319),
    @(45627.99999999999, 72.0, -57.45951644006173, 197.8488234274854),
    @(45641.99999999999, 72.0, -52.13826438032333, 196.5263487423765),
    @(45648.99999999999, 72.0, -60.28311188353717, 194.1086348966857),
    @(45655.99999999999, 72.0, -59.29776019589153, 199.6579538107232),
    @(45662.99999999999, 72.0, -56.59558764551767, 208.644250157083),
    @(45669.99999999999, 72.0, -50.87790940057959, 204.2945047443312),
    @(45676.99999999999, 72.0, -55.78225440617986, 193.7637317640966),
    @(45683.99999999999, 72.0, -59.58616225391219, 191.1789776097858),
    @(45690.99999999999, 72.0, -50.55083420477186, 194.1026264207607),
    @(45697.99999999999, 72.0, -51.01010537953076, 194.7462181324152)
)

$row = 2
foreach ($rec in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $rec[0]
    $wsForecast.Cells.Item($row, 2).Value = $rec[1]
    $wsForecast.Cells.Item($row, 3).Value = $rec[2]
    $wsForecast.Cells.Item($row, 4).Value = $rec[3]
    $row++
}
